$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text number format on D (price) and E (volume %) columns so that
# these numeric-looking / percent-looking strings are stored as literal text,
# matching the original inlineStr cells (not converted to numbers/percentages).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.78"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.34%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "29.76"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.31%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.290"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.23%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05749"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.74%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.650"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.72%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.232"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "6.30%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8574"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.17%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8562"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-2.06%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1380"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.00%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07089"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.06%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03220"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "12.54%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09353"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.35%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001528"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.17%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0005986"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.05%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005990"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.52%"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.06%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.194"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.36%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3157"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.38%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.03339"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2.80%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1298"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.20%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.478"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.12%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.1409"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2.08%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04119"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.89%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001222"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.33%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004170"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-18.04%"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.87%"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001448"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "-25.30%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03765"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.36%"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.16%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009193"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-5.79%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005276"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "3.16%"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.05%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.08983"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "26.54%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002876"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "6.15%"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.05%"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.05%"

# Rows 42 and 43: Coin/Link swap (CEJI <-> KickToken) plus new Price/Volume values
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.002393"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-5.03%"

$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002948"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-14.91%"
